$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 19608624
$ws.Range("I19").Value = 41667110
$ws.Range("K19").Value = 41667110
$ws.Range("M19").Value = -41666935
$ws.Range("H98").Value = 1433.25
$ws.Range("I98").Value = 983.25
$ws.Range("K98").Value = 983.25
$ws.Range("M98").Value = 514.75
$ws.Range("H100").Value = 1264.8286
$ws.Range("I100").Value = 1117.3462
$ws.Range("J100").Value = 1690.8889
$ws.Range("K100").Value = 1117.3462
$ws.Range("L100").Value = 1690.8889
$ws.Range("M100").Value = -576.3462
$ws.Range("N100").Value = -2772.8889
$ws.Range("H103").Value = 882.7895
$ws.Range("I103").Value = 600
$ws.Range("J103").Value = 983.7857
$ws.Range("K103").Value = 1800
$ws.Range("L103").Value = 2951.3571
$ws.Range("M103").Value = -1214
$ws.Range("N103").Value = -4123.3571
$ws.Range("H105").Value = 34871
$ws.Range("J105").Value = 34871
$ws.Range("L105").Value = 34871
$ws.Range("N105").Value = -41859
$ws.Range("H106").Value = 2756.8572
$ws.Range("I106").Value = 2756.8572
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2756.8572
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2125.8572
$ws.Range("N106").Value = ""
$ws.Range("H107").Value = 470.74194
$ws.Range("I107").Value = 361.875
$ws.Range("K107").Value = 361.875
$ws.Range("M107").Value = 1558.125
$ws.Range("H116").Value = 2085.25
$ws.Range("I116").Value = 1889.7
$ws.Range("J116").Value = 2280.8
$ws.Range("K116").Value = 1889.7
$ws.Range("L116").Value = 2280.8
$ws.Range("M116").Value = 1552.3
$ws.Range("N116").Value = -9164.799999999999
$ws.Range("H122").Value = 1433.25
$ws.Range("I122").Value = 983.25
$ws.Range("K122").Value = 2949.75
$ws.Range("M122").Value = -499.75

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1602.9524
$ws.Range("I45").Value = 1619.0526
$ws.Range("J45").Value = 1450
$ws.Range("K45").Value = 1619.0526
$ws.Range("L45").Value = 1450
$ws.Range("M45").Value = -1242.0526
$ws.Range("N45").Value = -2204
$ws.Range("H61").Value = 6199.8623
$ws.Range("I61").Value = 5035.32
$ws.Range("K61").Value = 5035.32
$ws.Range("M61").Value = -4823.32
$ws.Range("H102").Value = 2384.4546
$ws.Range("I102").Value = 2270
$ws.Range("K102").Value = 2270
$ws.Range("M102").Value = -648
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""
$ws.Range("H136").Value = 6199.8623
$ws.Range("I136").Value = 5035.32
$ws.Range("K136").Value = 15105.96
$ws.Range("M136").Value = -12555.96
$ws.Range("H137").Value = 40000
$ws.Range("J137").Value = 40000
$ws.Range("L137").Value = 40000
$ws.Range("N137").Value = -50200

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1369.3077
$ws.Range("I134").Value = 1401
$ws.Range("J134").Value = 989
$ws.Range("K134").Value = 4203
$ws.Range("L134").Value = 2967
$ws.Range("M134").Value = -1668
$ws.Range("N134").Value = -8037

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2197.111
$ws.Range("I16").Value = 1700.9
$ws.Range("K16").Value = 1700.9
$ws.Range("M16").Value = -1413.9
$ws.Range("H31").Value = 2069.9722
$ws.Range("I31").Value = 1311.6451
$ws.Range("J31").Value = 6771.6
$ws.Range("K31").Value = 1311.6451
$ws.Range("L31").Value = 6771.6
$ws.Range("M31").Value = -1016.6451
$ws.Range("N31").Value = -7361.6
$ws.Range("H34").Value = 2069.9722
$ws.Range("I34").Value = 1311.6451
$ws.Range("J34").Value = 6771.6
$ws.Range("K34").Value = 1311.6451
$ws.Range("L34").Value = 6771.6
$ws.Range("M34").Value = -1109.6451
$ws.Range("N34").Value = -7175.6
$ws.Range("H99").Value = 1680.2142
$ws.Range("I99").Value = 1377.3
$ws.Range("J99").Value = 2437.5
$ws.Range("K99").Value = 1377.3
$ws.Range("L99").Value = 2437.5
$ws.Range("M99").Value = 120.7
$ws.Range("N99").Value = -5433.5
$ws.Range("H105").Value = 680.8182
$ws.Range("I105").Value = 563.75
$ws.Range("J105").Value = 993
$ws.Range("K105").Value = 563.75
$ws.Range("L105").Value = 993
$ws.Range("M105").Value = 1183.25
$ws.Range("N105").Value = -4487
$ws.Range("H107").Value = 604.8261
$ws.Range("I107").Value = 347.2857
$ws.Range("J107").Value = 717.5
$ws.Range("K107").Value = 347.2857
$ws.Range("L107").Value = 717.5
$ws.Range("M107").Value = 1572.7143
$ws.Range("N107").Value = -4557.5
$ws.Range("H113").Value = 2197.111
$ws.Range("I113").Value = 1700.9
$ws.Range("K113").Value = 1700.9
$ws.Range("M113").Value = 469.0999999999999
$ws.Range("H122").Value = 5766.2856
$ws.Range("I122").Value = 5979.077
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 17937.231
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -15487.231
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 1680.2142
$ws.Range("I126").Value = 1377.3
$ws.Range("J126").Value = 2437.5
$ws.Range("K126").Value = 4131.9
$ws.Range("L126").Value = 7312.5
$ws.Range("M126").Value = -1661.9
$ws.Range("N126").Value = -12252.5
$ws.Range("H138").Value = 36059.168
$ws.Range("J138").Value = 36059.168
$ws.Range("L138").Value = 36059.168
$ws.Range("N138").Value = -46339.168

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 221.5
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H94").Value = 15000
$ws.Range("J94").Value = 15000
$ws.Range("L94").Value = 45000
$ws.Range("N94").Value = -46352
$ws.Range("H95").Value = 6006
$ws.Range("J95").Value = 6000
$ws.Range("L95").Value = 18000
$ws.Range("N95").Value = -22118
$ws.Range("H97").Value = 1530.7
$ws.Range("I97").Value = 501
$ws.Range("J97").Value = 3933.3333
$ws.Range("K97").Value = 1503
$ws.Range("L97").Value = 11799.9999
$ws.Range("M97").Value = -1007
$ws.Range("N97").Value = -12791.9999
$ws.Range("H98").Value = 459.95
$ws.Range("I98").Value = 400
$ws.Range("J98").Value = 999.5
$ws.Range("K98").Value = 1200
$ws.Range("L98").Value = 2998.5
$ws.Range("M98").Value = 298
$ws.Range("N98").Value = -5994.5
$ws.Range("H100").Value = 3938
$ws.Range("J100").Value = 3938
$ws.Range("L100").Value = 11814
$ws.Range("N100").Value = -13436
$ws.Range("H101").Value = 7219.778
$ws.Range("J101").Value = 7219.778
$ws.Range("L101").Value = 21659.334
$ws.Range("N101").Value = -26527.334
$ws.Range("H103").Value = 805.75
$ws.Range("I103").Value = 408
$ws.Range("J103").Value = 1999
$ws.Range("K103").Value = 1224
$ws.Range("L103").Value = 5997
$ws.Range("M103").Value = -345
$ws.Range("N103").Value = -7755
$ws.Range("H104").Value = 2640.2222
$ws.Range("I104").Value = 440
$ws.Range("J104").Value = 3268.8572
$ws.Range("K104").Value = 1320
$ws.Range("L104").Value = 9806.571599999999
$ws.Range("M104").Value = 1301
$ws.Range("N104").Value = -15048.5716

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 6000
$ws.Range("I92").Value = 6000
$ws.Range("K92").Value = 6000
$ws.Range("M92").Value = -4128
$ws.Range("H94").Value = 30400
$ws.Range("J94").Value = 30400
$ws.Range("L94").Value = 30400
$ws.Range("N94").Value = -31752
$ws.Range("H97").Value = 2472
$ws.Range("I97").Value = 2472
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2472
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1976
$ws.Range("N97").Value = ""
$ws.Range("H102").Value = 3679.7896
$ws.Range("I102").Value = 2847.7334
$ws.Range("K102").Value = 2847.7334
$ws.Range("M102").Value = -1225.7334
$ws.Range("H122").Value = 4748.3335
$ws.Range("I122").Value = 5025
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 15075
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -12625
$ws.Range("N122").Value = -13750

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2762.4092
$ws.Range("I40").Value = 2489
$ws.Range("J40").Value = 3692
$ws.Range("K40").Value = 2489
$ws.Range("L40").Value = 3692
$ws.Range("M40").Value = -2353
$ws.Range("N40").Value = -3964
$ws.Range("H122").Value = 5995.5713
$ws.Range("I122").Value = 5474.353
$ws.Range("K122").Value = 16423.059
$ws.Range("M122").Value = -13973.059
$ws.Range("H125").Value = 58000
$ws.Range("J125").Value = 58000
$ws.Range("L125").Value = 58000
$ws.Range("N125").Value = -67840
$ws.Range("H134").Value = 16600
$ws.Range("I134").Value = 15000
$ws.Range("J134").Value = 19800
$ws.Range("K134").Value = 15000
$ws.Range("L134").Value = 19800
$ws.Range("M134").Value = -9930
$ws.Range("N134").Value = -29940
$ws.Range("H136").Value = 4401.222
$ws.Range("I136").Value = 2620.4583
$ws.Range("K136").Value = 7861.374899999999
$ws.Range("M136").Value = -5311.374899999999
$ws.Range("H141").Value = 49374
$ws.Range("J141").Value = 49374
$ws.Range("L141").Value = 49374
$ws.Range("N141").Value = -59734

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 35959
$ws.Range("J41").Value = 35959
$ws.Range("L41").Value = 35959
$ws.Range("N41").Value = -36739
$ws.Range("H122").Value = 758.8
$ws.Range("I122").Value = 723.2308
$ws.Range("K122").Value = 2169.6924
$ws.Range("M122").Value = 280.3076000000001
$ws.Range("H126").Value = 1355.9131
$ws.Range("I126").Value = 1326.6364
$ws.Range("K126").Value = 3979.9092
$ws.Range("M126").Value = -1509.9092
